$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MONTAGEM")
[void]$ws.Activate()

# Update the BOM description for R3 (row 22) and R4 (row 23) from the
# 330 ohm resistor to the 1 kOhm resistor per the 2.1 schematic update.
# Copy from B25 (already "Resistor 1 kOhm, 1/4 W - CR25") so the cells
# reuse the existing rich-text shared string instead of minting a new one.
[void]$ws.Range("B25").Copy($ws.Range("B22"))
[void]$ws.Range("B25").Copy($ws.Range("B23"))

# Match the author's final view state (scrolled/selected cell).
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B20").Select()
